# Insert a new "Force" row above the existing "Upload" row on the
# Property1 sheet (sheet1). The new row mirrors the layout/format of the
# row directly above it (the "Ref" row) and pushes "Upload", the "Desc"
# row, and the SqlServer_1 sample-data row down by one each.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 8 (the "Upload" row), shifting everything
# from row 8 downward (Upload, Desc row, sample-data row) down by one.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the "Force" entry (boolean flags default
# to FALSE, matching every other row of this block).
$ws.Cells.Item(8, 1).Value = "Force"
$ws.Cells.Item(8, 2).Value = $false
$ws.Cells.Item(8, 3).Value = $false
$ws.Cells.Item(8, 4).Value = $false
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(8, 6).Value = $false
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = $false
$ws.Cells.Item(8, 9).Value = $false

# Match formatting of the row above (row 7) for the new row.
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)  # xlPasteFormats

# (The data validation ranges "A6:A8" and "B6:I6 B7:J8" automatically
# grow to "A6:A9" and "B6:I6 B7:J9" as part of the row insert/shift
# above, since row 8 falls inside both validated ranges.)

# The freeze pane needs to move down one row to stay anchored below the
# header block (was frozen after row 9, now after row 10). Toggle it off
# and back on at the new anchor so Excel recomputes the split position.
$aw = $excel.ActiveWindow
$aw.FreezePanes = $false
$ws.Range("A11").Select()
$aw.FreezePanes = $true

# Restore the selection to A9, matching the saved view state.
$ws.Range("A9").Select()
